$d = $word.ActiveDocument

$replacements = @(
    @("235÷8=", "787÷3="),
    @("910÷9=", "307÷5="),
    @("822÷4=", "633÷5="),
    @("432÷2=", "169÷3="),
    @("304÷3=", "136÷6="),
    @("122÷3=", "830÷2="),
    @("940÷8=", "121÷3="),
    @("206÷4=", "303÷3="),
    @("203÷8=", "580÷4="),
    @("158÷2=", "990÷6="),
    @("773÷8=", "124÷6="),
    @("854÷2=", "682÷6="),
    @("724÷3=", "490÷3="),
    @("347÷6=", "987÷5="),
    @("990÷9=", "457÷3="),
    @("209÷2=", "990÷3="),
    @("954÷3=", "824÷2="),
    @("496÷7=", "695÷8="),
    @("204÷9=", "139÷3="),
    @("781÷2=", "298÷2="),
    @("693÷5=", "168÷8="),
    @("445÷6=", "350÷9="),
    @("992÷3=", "603÷3="),
    @("398÷8=", "829÷4="),
    @("295÷3=", "295÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
